# Replace the "develop/obtain/recruit/exchange/involve/collaborate" + "PC_ENT"
# indicator/unit codes in columns B & C (rows 2-19) with the new
# "E_SM_*" social-media indicator codes and "PC_ENT_SM" unit code.
# Columns D (sizen_r2) and E (geo=EU28) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
  2 = @("E_SM_PADVERT", "PC_ENT_SM", "L_C10_S951_XK", "EU28")
  3 = @("E_SM_PADVERT", "PC_ENT_SM", "M_C10_S951_XK", "EU28")
  4 = @("E_SM_PADVERT", "PC_ENT_SM", "S_C10_S951_XK", "EU28")
  5 = @("E_SM_PCUQOR", "PC_ENT_SM", "L_C10_S951_XK", "EU28")
  6 = @("E_SM_PCUQOR", "PC_ENT_SM", "M_C10_S951_XK", "EU28")
  7 = @("E_SM_PCUQOR", "PC_ENT_SM", "S_C10_S951_XK", "EU28")
  8 = @("E_SM_PRCR", "PC_ENT_SM", "L_C10_S951_XK", "EU28")
  9 = @("E_SM_PRCR", "PC_ENT_SM", "M_C10_S951_XK", "EU28")
  10 = @("E_SM_PRCR", "PC_ENT_SM", "S_C10_S951_XK", "EU28")
  11 = @("E_SM_PEXCHVOC", "PC_ENT_SM", "L_C10_S951_XK", "EU28")
  12 = @("E_SM_PEXCHVOC", "PC_ENT_SM", "M_C10_S951_XK", "EU28")
  13 = @("E_SM_PEXCHVOC", "PC_ENT_SM", "S_C10_S951_XK", "EU28")
  14 = @("E_SM_PCUDEV", "PC_ENT_SM", "L_C10_S951_XK", "EU28")
  15 = @("E_SM_PCUDEV", "PC_ENT_SM", "M_C10_S951_XK", "EU28")
  16 = @("E_SM_PCUDEV", "PC_ENT_SM", "S_C10_S951_XK", "EU28")
  17 = @("E_SM_PBPCOLL", "PC_ENT_SM", "L_C10_S951_XK", "EU28")
  18 = @("E_SM_PBPCOLL", "PC_ENT_SM", "M_C10_S951_XK", "EU28")
  19 = @("E_SM_PBPCOLL", "PC_ENT_SM", "S_C10_S951_XK", "EU28")
}

# Write column C (new unit "PC_ENT_SM") before column B on each row so the
# shared-strings table is populated in the same order as the target file
# (PC_ENT_SM ends up right after S_C10_S951_XK, before the E_SM_* indicators).
foreach ($row in ($map.Keys | Sort-Object)) {
  $vals = $map[$row]
  $ws.Cells.Item($row, 3).Value = $vals[1]
  $ws.Cells.Item($row, 2).Value = $vals[0]
  $ws.Cells.Item($row, 4).Value = $vals[2]
  $ws.Cells.Item($row, 5).Value = $vals[3]
}

# Columns B & C now hold longer strings -> widen them (mirrors the author
# re-running "AutoFit column width" after the content change).
$ws.Columns.Item(2).ColumnWidth = 14.6
$ws.Columns.Item(3).ColumnWidth = 10.0

# Selection moved from M9 to I11 (also clears the stale topLeftCell scroll state).
$ws.Range("I11").Select() | Out-Null
